# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - copy the format of the existing header cell (AC1)
# so the new header cells share the same bold/centered/bordered style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-54: Wins=92, Losses=71, Ties=0
for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 30).Value = 92  # AD
    $ws.Cells.Item($r, 31).Value = 71  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
